$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for the new Wins/Losses/Ties columns, copying the
# formatting from the existing header row (e.g. AC1) so the new headers
# match the rest of row 1 (bold, bordered, centered).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-48) gets the team's record for the season.
$ws.Range("AD2:AD48").Value = 86
$ws.Range("AE2:AE48").Value = 76
$ws.Range("AF2:AF48").Value = 0
